# Fix the "2050" / "2041-2050" column-header label (which had been
# accidentally overwritten with a stray numeric value) and drop the
# "Total" summary row from the scenario tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets

# --- helper: write `text` into `cell` as a genuine text value (never an
# auto-coerced number) while leaving the cell's existing style untouched.
# Writing Range.Value directly with a numeric-looking string (e.g. "2050")
# gets auto-converted to a number by the COM layer, and forcing text via
# NumberFormat="@" mints a brand-new cell style (bumping the cellXfs
# index). Going through a scratch cell + Copy/PasteSpecial(xlPasteValues)
# carries over only the *value* (already typed as text), so the
# destination cell's style index is left exactly as it was.
function Set-TextValue {
    param($sheet, $row, $col, [string]$text)

    $scratch = $sheet.Cells.Item(5000, 50)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()

    $target = $sheet.Cells.Item($row, $col)
    $target.PasteSpecial(-4163)  # xlPasteValues

    $scratch.Clear()
}

# Sheets 1-3: "...,2015,2030,2040,<broken>" -> relabel E1 "2050"
foreach ($idx in 1,2,3) {
    $sheet = $ws.Item($idx)
    Set-TextValue $sheet 1 5 "2050"
    $sheet.Rows.Item(13).Delete()
}

# Sheet 4: "...,2015-2030,2031-2040,<broken>" -> relabel E1 "2041-2050"
$sheet4 = $ws.Item(4)
Set-TextValue $sheet4 1 5 "2041-2050"
$sheet4.Rows.Item(13).Delete()

# Sheet 5: only the mislabeled header needs fixing (no Total row here)
$sheet5 = $ws.Item(5)
Set-TextValue $sheet5 1 5 "2050"

# Sheet 6: no E1 header cell (only column B) - just drop the Total row
$sheet6 = $ws.Item(6)
$sheet6.Rows.Item(4).Delete()
